$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 6).Value = 4.7  # F2: 4.8 -> 4.7
$ws.Cells.Item(2, 8).Value = 1.51  # H2: 1.46 -> 1.51
$ws.Cells.Item(2, 11).Value = 8  # K2: 9.199999999999999 -> 8
$ws.Cells.Item(2, 15).Value = 1.17  # O2: 1.19 -> 1.17
$ws.Cells.Item(2, 17).Value = 1.54  # Q2: 1.51 -> 1.54
$ws.Cells.Item(2, 19).Value = 2.3  # S2: 2.26 -> 2.3
$ws.Cells.Item(2, 23).Value = 1.07  # W2: 1.04 -> 1.07
$ws.Cells.Item(3, 8).Value = 1.05  # H3: 1.04 -> 1.05
$ws.Cells.Item(3, 18).Value = 1.18  # R3: 1.2 -> 1.18
$ws.Cells.Item(4, 6).Value = 1.34  # F4: 1.35 -> 1.34
$ws.Cells.Item(4, 8).Value = 11  # H4: 10.5 -> 11
$ws.Cells.Item(4, 9).Value = 12  # I4: 11.5 -> 12
$ws.Cells.Item(4, 20).Value = 1.91  # T4: 1.92 -> 1.91
$ws.Cells.Item(4, 25).Value = 1000  # Y4: 42 -> 1000
$ws.Cells.Item(4, 26).Value = 110  # Z4: 100 -> 110
$ws.Cells.Item(4, 31).Value = 160  # AE4: 150 -> 160
$ws.Cells.Item(4, 35).Value = 130  # AI4: 120 -> 130
$ws.Cells.Item(5, 6).Value = 5.2  # F5: 5.3 -> 5.2
$ws.Cells.Item(5, 9).Value = 1.76  # I5: 1.75 -> 1.76
$ws.Cells.Item(5, 10).Value = 4.2  # J5: 4.1 -> 4.2
$ws.Cells.Item(5, 18).Value = 1.54  # R5: 1.55 -> 1.54
$ws.Cells.Item(5, 20).Value = 1.73  # T5: 1.72 -> 1.73
$ws.Cells.Item(5, 23).Value = 1.23  # W5: 1.22 -> 1.23
$ws.Cells.Item(7, 12).Value = 1.26  # L7: 1.3 -> 1.26
$ws.Cells.Item(8, 7).Value = 2.4  # G8: 2.42 -> 2.4
$ws.Cells.Item(8, 9).Value = 3.3  # I8: 3.35 -> 3.3
$ws.Cells.Item(8, 14).Value = 4.2  # N8: 4.1 -> 4.2
$ws.Cells.Item(8, 23).Value = 1.71  # W8: 1.7 -> 1.71
$ws.Cells.Item(8, 24).Value = 16.5  # X8: 15 -> 16.5
$ws.Cells.Item(8, 29).Value = 8  # AC8: 7.8 -> 8
$ws.Cells.Item(9, 18).Value = 2.06  # R9: 2.08 -> 2.06
$ws.Cells.Item(9, 22).Value = 1.1  # V9: 1.11 -> 1.1
$ws.Cells.Item(9, 23).Value = 3.9  # W9: 3.85 -> 3.9
$ws.Cells.Item(9, 25).Value = 60  # Y9: 55 -> 60
$ws.Cells.Item(10, 8).Value = 9.6  # H10: 10 -> 9.6
$ws.Cells.Item(10, 10).Value = 5.4  # J10: 5.5 -> 5.4
$ws.Cells.Item(10, 29).Value = 13  # AC10: 12.5 -> 13
$ws.Cells.Item(10, 33).Value = 10.5  # AG10: 10 -> 10.5
$ws.Cells.Item(10, 34).Value = 23  # AH10: 24 -> 23
$ws.Cells.Item(11, 15).Value = 1.19  # O11: 1.2 -> 1.19
$ws.Cells.Item(11, 36).Value = 18  # AJ11: 17.5 -> 18
$ws.Cells.Item(12, 14).Value = 4.5  # N12: 4.6 -> 4.5
$ws.Cells.Item(12, 19).Value = 2.96  # S12: 2.98 -> 2.96
$ws.Cells.Item(12, 21).Value = 2.4  # U12: 2.42 -> 2.4
$ws.Cells.Item(13, 7).Value = 8.199999999999999  # G13: 8 -> 8.199999999999999
$ws.Cells.Item(13, 10).Value = 5.4  # J13: 5.3 -> 5.4
$ws.Cells.Item(13, 16).Value = 2.74  # P13: 2.72 -> 2.74
$ws.Cells.Item(13, 17).Value = 1.54  # Q13: 1.55 -> 1.54
$ws.Cells.Item(13, 37).Value = 90  # AK13: 95 -> 90
$ws.Cells.Item(13, 39).Value = 90  # AM13: 95 -> 90
$ws.Cells.Item(13, 41).Value = 5  # AO13: 4.9 -> 5
$ws.Cells.Item(14, 16).Value = 2.4  # P14: 2.38 -> 2.4
$ws.Cells.Item(16, 21).Value = 2.22  # U16: 2.2 -> 2.22
$ws.Cells.Item(16, 25).Value = 980  # Y16: 970 -> 980
$ws.Cells.Item(16, 26).Value = 980  # Z16: 30 -> 980
$ws.Cells.Item(16, 30).Value = 980  # AD16: 970 -> 980
$ws.Cells.Item(16, 31).Value = 980  # AE16: 44 -> 980
$ws.Cells.Item(16, 32).Value = 980  # AF16: 970 -> 980
$ws.Cells.Item(16, 34).Value = 980  # AH16: 970 -> 980
$ws.Cells.Item(16, 35).Value = 980  # AI16: 50 -> 980
$ws.Cells.Item(16, 36).Value = 980  # AJ16: 25 -> 980
$ws.Cells.Item(16, 37).Value = 980  # AK16: 21 -> 980
$ws.Cells.Item(16, 38).Value = 980  # AL16: 34 -> 980
$ws.Cells.Item(16, 41).Value = 980  # AO16: 42 -> 980
$ws.Cells.Item(18, 6).Value = 1.73  # F18: 1.81 -> 1.73
$ws.Cells.Item(18, 8).Value = 4.5  # H18: 3.85 -> 4.5
$ws.Cells.Item(18, 9).Value = 7.6  # I18: 7.8 -> 7.6
$ws.Cells.Item(18, 10).Value = 1.2  # J18: 2.44 -> 1.2
$ws.Cells.Item(18, 11).Value = 5.1  # K18: 500 -> 5.1
$ws.Cells.Item(18, 18).Value = 1.22  # R18: 1.2 -> 1.22
$ws.Cells.Item(18, 22).Value = 1.19  # V18: 1.18 -> 1.19
$ws.Cells.Item(19, 6).Value = 2.04  # F19: 2.1 -> 2.04
$ws.Cells.Item(19, 23).Value = 1.8  # W19: 1.78 -> 1.8
$ws.Cells.Item(19, 32).Value = 12  # AF19: 980 -> 12
$ws.Cells.Item(19, 37).Value = 980  # AK19: 40 -> 980
$ws.Cells.Item(20, 15).Value = 1.31  # O20: 1.3 -> 1.31
$ws.Cells.Item(20, 16).Value = 1.96  # P20: 1.95 -> 1.96
$ws.Cells.Item(20, 26).Value = 15  # Z20: 1000 -> 15
$ws.Cells.Item(20, 35).Value = 1000  # AI20: 38 -> 1000
$ws.Cells.Item(20, 41).Value = 19  # AO20: 1000 -> 19
$ws.Cells.Item(21, 16).Value = 2.12  # P21: 2.14 -> 2.12
$ws.Cells.Item(21, 19).Value = 2.6  # S21: 2.54 -> 2.6
$ws.Cells.Item(21, 20).Value = 1.81  # T21: 1.8 -> 1.81
$ws.Cells.Item(21, 22).Value = 1.15  # V21: 1.13 -> 1.15
$ws.Cells.Item(21, 29).Value = 12.5  # AC21: 1000 -> 12.5
$ws.Cells.Item(21, 33).Value = 12.5  # AG21: 1000 -> 12.5
